$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "26.133.40"
$ws.Range("D3").Value = "1.658.12"
$ws.Range("D4").Formula = "=""1.003"""
$ws.Range("D5").Formula = "=""216.57"""
$ws.Range("D6").Formula = "=""0.5169"""
$ws.Range("D8").Formula = "=""0.2636"""
$ws.Range("D9").Formula = "=""0.06273"""
$ws.Range("D10").Formula = "=""20.76"""
$ws.Range("D11").Formula = "=""0.07714"""
$ws.Range("D12").Value = "1.652.92"
$ws.Range("D13").Formula = "=""4.428"""
$ws.Range("D14").Value = "1.884.67"
$ws.Range("D15").Formula = "=""0.5416"""
$ws.Range("D16").Value = "0.0₅8112"
$ws.Range("D17").Formula = "=""64.80"""
$ws.Range("D18").Value = "26.155.76"
$ws.Range("D20").Formula = "=""4.619"""
$ws.Range("D21").Formula = "=""191.62"""
$ws.Range("D23").Formula = "=""6.014"""
$ws.Range("D25").Formula = "=""139.78"""
$ws.Range("D26").Formula = "=""0.1225"""
$ws.Range("D27").Formula = "=""7.188"""
$ws.Range("D28").Formula = "=""16.08"""
$ws.Range("D29").Formula = "=""1.406"""
$ws.Range("D30").Formula = "=""0.05964"""
$ws.Range("D32").Formula = "=""3.554"""
$ws.Range("D33").Formula = "=""3.255"""
$ws.Range("D34").Formula = "=""1.601"""
$ws.Range("D35").Formula = "=""0.9644"""
$ws.Range("D36").Formula = "=""2.425"""
$ws.Range("D37").Formula = "=""2.769"""
$ws.Range("D38").Formula = "=""0.5687"""
$ws.Range("D39").Formula = "=""0.01590"""
$ws.Range("D40").Formula = "=""5.954"""
$ws.Range("D41").Formula = "=""0.8542"""
$ws.Range("D43").Value = "1.006.13"
$ws.Range("D44").Formula = "=""100.46"""
$ws.Range("D45").Value = "1.799.21"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("D47").Formula = "=""56.71"""
$ws.Range("D49").Formula = "=""7.993"""
$ws.Range("D51").Formula = "=""1.446"""

# Convert formula-based text entries back to literal text values
# (keeps default style, avoids Excel auto-coercing numeric-looking strings to numbers)
$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E10").Value = "  -4.84%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("E23").Value = "  -4.94%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("E30").Value = "  -5.26%  "
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("E33").Value = "  -5.48%  "
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("E35").Value = "  -4.92%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  -7.84%  "
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("E51").Value = "  -3.99%  "
